$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1920478055701966
$ws.Range("I2").Value = 0.2931153080917904
$ws.Range("J2").Value = 0.2240828618179215
$ws.Range("K2").Value = 0.7639291549634947
$ws.Range("O2").Value = 0.3018189596204621
$ws.Range("G3").Value = 0.8440802290461885
$ws.Range("I3").Value = 0.9486587492742637
$ws.Range("J3").Value = 0.8030715638437396
$ws.Range("K3").Value = 0.9200936938865002
$ws.Range("M3").Value = 1.344919739557854
$ws.Range("O3").Value = 1.147772757933678
$ws.Range("G4").Value = 0.3205563771652512
$ws.Range("I4").Value = 0.2765914636051733
$ws.Range("J4").Value = 0.4265827512322169
$ws.Range("K4").Value = 0.706912652919841
$ws.Range("M4").Value = 0.2898130628153743
$ws.Range("O4").Value = 0.6361204696781061
$ws.Range("G5").Value = 1.438599549706423
$ws.Range("I5").Value = 1.227672168773977
$ws.Range("J5").Value = 1.902151601688216
$ws.Range("K5").Value = 0.8355873218281257
$ws.Range("M5").Value = 0.7818175215325235
$ws.Range("O5").Value = 1.164045658562616
$ws.Range("G6").Value = 0.9131169179845734
$ws.Range("I6").Value = 0.6927983464168281
$ws.Range("J6").Value = 0.499637322800165
$ws.Range("K6").Value = 0.6561596113299912
$ws.Range("M6").Value = 0.4073430323889371
$ws.Range("O6").Value = 0.3617692475965423
$ws.Range("G7").Value = 0.2339345087460386
$ws.Range("I7").Value = 0.2340170282194848
$ws.Range("J7").Value = 0.2211828025353424
$ws.Range("K7").Value = 0.7255852489729558
$ws.Range("M7").Value = 0.3388181849359906
$ws.Range("O7").Value = 0.4219135222636656
$ws.Range("G8").Value = 0.5218584798125326
$ws.Range("I8").Value = 0.4269770201465697
$ws.Range("J8").Value = 0.3478409681710704
$ws.Range("K8").Value = 0.6713204581193503
$ws.Range("M8").Value = 0.2777087986815598
$ws.Range("O8").Value = 0.342433425251829
$ws.Range("G9").Value = 0.7166187989284143
$ws.Range("I9").Value = 0.5093922466907995
$ws.Range("J9").Value = 0.3622002300607757
$ws.Range("K9").Value = 0.6506051046548891
$ws.Range("M9").Value = 0.3202129035300945
$ws.Range("O9").Value = 0.4218013207848155
$ws.Range("G10").Value = 0.6024356727075499
$ws.Range("I10").Value = 0.7197157884090679
$ws.Range("J10").Value = 1.156404960599048
$ws.Range("K10").Value = 0.8599616982087048
$ws.Range("M10").Value = 1.052450656498203
$ws.Range("O10").Value = 0.7045524970028708
$ws.Range("G11").Value = 1.862080997882652
$ws.Range("I11").Value = 1.584663472020739
$ws.Range("J11").Value = 1.674600097977642
$ws.Range("K11").Value = 1.146178872540844
$ws.Range("M11").Value = 1.20733943113395
$ws.Range("O11").Value = 1.880874184483762
$ws.Range("G12").Value = 2.439748241263175
$ws.Range("I12").Value = 2.156943198439215
$ws.Range("J12").Value = 1.086398644454755
$ws.Range("K12").Value = 1.72620320147629
$ws.Range("M12").Value = 1.801475110378239
$ws.Range("O12").Value = 1.419224967255788
$ws.Range("G13").Value = 2.209278834976443
$ws.Range("I13").Value = 1.941628516358671
$ws.Range("J13").Value = 0.8538861526896857
$ws.Range("K13").Value = 1.450516649476119
$ws.Range("M13").Value = 1.565679653725108
$ws.Range("O13").Value = 1.348002749600517
$ws.Range("G14").Value = 0.5594845440567042
$ws.Range("I14").Value = 0.6066304995240968
$ws.Range("J14").Value = 0.5877663687395481
$ws.Range("K14").Value = 0.6295399261422575
$ws.Range("O14").Value = 0.4781814993678736
$ws.Range("G15").Value = 0.2330449601877848
$ws.Range("I15").Value = 0.3116334807214568
$ws.Range("J15").Value = 0.3074736838596404
$ws.Range("K15").Value = 0.6770930661394439
$ws.Range("M15").Value = 0.4158797868320672
$ws.Range("O15").Value = 0.2546694730609776
$ws.Range("G16").Value = 0.8301659967663215
$ws.Range("I16").Value = 0.6487811912476636
$ws.Range("J16").Value = 0.6940797106820545
$ws.Range("K16").Value = 0.8630725898716562
$ws.Range("M16").Value = 0.5047890299838209
$ws.Range("O16").Value = 0.6396820294032204
$ws.Range("G17").Value = 0.2217847782587894
$ws.Range("I17").Value = 0.2636546328138152
$ws.Range("J17").Value = 0.2897748490341038
$ws.Range("K17").Value = 0.7429760596458208
$ws.Range("M17").Value = 0.2408162117359516
$ws.Range("O17").Value = 0.2508553772214462
$ws.Range("G18").Value = 0.1914774984736885
$ws.Range("I18").Value = 0.2688765158883483
$ws.Range("J18").Value = 0.2682123033864766
$ws.Range("K18").Value = 0.6990447305795762
$ws.Range("M18").Value = 0.3146403115582309
$ws.Range("O18").Value = 0.2349082125301026
$ws.Range("G19").Value = 0.347339177607937
$ws.Range("I19").Value = 0.4181706787803032
$ws.Range("J19").Value = 0.3914062694776969
$ws.Range("K19").Value = 0.654271755691425
$ws.Range("M19").Value = 0.5703237083467829
$ws.Range("O19").Value = 0.3521912712262327
$ws.Range("G20").Value = 0.2094056339755848
$ws.Range("I20").Value = 0.2354635038536856
$ws.Range("J20").Value = 0.2766604675248611
$ws.Range("K20").Value = 0.7388328882395838
$ws.Range("M20").Value = 0.2410705361697636
$ws.Range("O20").Value = 0.2784748908302692
$ws.Range("G21").Value = 0.2803610324661884
$ws.Range("I21").Value = 0.2862780413145777
$ws.Range("J21").Value = 0.3010800891458826
$ws.Range("K21").Value = 0.765275624254403
$ws.Range("M21").Value = 0.2558664925255891
$ws.Range("O21").Value = 0.2579204342093958
$ws.Range("G22").Value = 0.5121723496042556
$ws.Range("I22").Value = 0.4356322078536043
$ws.Range("J22").Value = 0.4601860729658694
$ws.Range("K22").Value = 0.8114704239280416
$ws.Range("M22").Value = 0.3451142342961169
$ws.Range("O22").Value = 0.4130297764902056
$ws.Range("G23").Value = 0.4364556868657948
$ws.Range("I23").Value = 0.3493979849081216
$ws.Range("J23").Value = 0.3934539751368294
$ws.Range("K23").Value = 0.7988085137592207
$ws.Range("M23").Value = 0.3138012994850925
$ws.Range("O23").Value = 0.3787537548649644
$ws.Range("G24").Value = 0.5717124853322739
$ws.Range("I24").Value = 0.6146609287772477
$ws.Range("J24").Value = 0.6155464648401392
$ws.Range("K24").Value = 0.635463278372507
$ws.Range("M24").Value = 0.7963150737309872
$ws.Range("O24").Value = 0.5675910868027001
$ws.Range("G25").Value = 0.2665965921611493
$ws.Range("I25").Value = 0.2853736705985076
$ws.Range("J25").Value = 0.317474040658008
$ws.Range("K25").Value = 0.7631565096809071
$ws.Range("M25").Value = 0.3574792898582413
$ws.Range("O25").Value = 0.3574792898582413
